# Edit: split the leading "MiU 1303 - 03/" run so "MiU" is wrapped in
# proofErr spellStart/spellEnd marks (as Word's spell-checker would do
# when it flags "MiU" as a misspelling), and bump the day-of-month run
# from "21" to "28". The trailing "_GoBack" bookmark and the "/13" run
# that follow must be preserved untouched.

$d = $word.ActiveDocument

$enDash = [char]0x2013

# Locate the paragraph that holds the "MiU 1303 ... " course/date line
# instead of a hard-coded index, so the script keeps working even if
# earlier content in the document shifts paragraph numbers around.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("MiU")) {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'MiU ...' paragraph"
}

$range = $target.Range

# Replace the whole paragraph's contents with the restructured run/proofErr
# layout. InsertXML replaces everything in the given range, so the
# bookmark and the trailing "/13" run are re-supplied verbatim here.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
  'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + `
  'w14:paraId="14BC5F0E" w14:textId="52849403" w:rsidR="00BF56E8" w:rsidRDefault="00314FA3">' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>MiU</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> 1303 ' + $enDash + ' 03/</w:t></w:r>' + `
  '<w:r w:rsidR="00EB6827"><w:t>28</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '<w:r w:rsidR="00BF56E8"><w:t>/13</w:t></w:r>' + `
  '</w:p>'

$range.InsertXML($xml)
